# Add scenario 9 - Lomas CET + Demand + Drug budget
$wb = $excel.ActiveWorkbook

# --- Sheet 1 (output_and_resourceuse): append summary row 9 ---
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A9").Value = 'CET ($164.7) + Demand constraint + Drug budget'
$ws1.Range("B9").Value = 89
$ws1.Range("C9").Value = 72
$ws1.Range("D9").Value = 67918104.11
$ws1.Range("E9").Value = 123.64
$ws1.Range("F9").Value = 1
$ws1.Range("G9").Value = 2.98
$ws1.Range("H9").Value = 1.52
$ws1.Range("I9").Value = 4.56
$ws1.Range("J9").Value = 0

# --- Sheet 2 (optimal_coverage): append column K for the new scenario ---
$ws2 = $wb.Worksheets.Item(2)

# Header cell K1 - styled like the rest of row 1 (bold + centered)
$ws2.Range("K1").Value = 'CET ($164.7) + Demand constraint + Drug budget'
$ws2.Range("K1").Font.Bold = $ws2.Range("J1").Font.Bold
$ws2.Range("K1").HorizontalAlignment = $ws2.Range("J1").HorizontalAlignment

# Per-intervention coverage values for the new scenario (rows 2-142)
$kValues = @(0.4699999999999921,0,0,0.7100000000004439,0,0.5999999999998817,0.9999999999993293,0.4699999999997057,0.6999999999999998,1.000000000000037,0,0,0,0,0,0,0.4,0.6,0,0.6000000000000001,0.25,0,0,1,0.8999999999999999,0,0,0,0.7,0,0,0.9000000000000001,0.6000000000001343,0.6000000000002146,0.5999999999999478,0.899999999999381,0,0,0.799999999999664,0.800000000000141,0.8000000000000136,0.9999999999994368,0.622495157119452,0.9499999999999997,0.95,0.9499999999999998,0.9499999999999997,0.5999999999998286,0.6000000000001128,0,0,0.5999999999996756,0.5999999999999488,0.6000000000001128,0.599999999999744,0.6000000000000002,0.5000000000003847,0.8000000000006154,0.05000000000000001,0.95,0.9500000000002184,0.9500000000000001,0,0.9499999999993444,0.9499999999994824,0.95,0.9500000000004576,0.9499999999999998,0.9500000000004774,0.9499999999997762,0,0,0,0.9499999999994906,0.95,0.95,0,0,0,0,0,0,0,0,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0.9999999999994607,0.5000000000002506,0,0,0.8999999999997979,0.5000000000001524,0.9399999999999999,0,0,0,0.9400000000000002,0.9399999999999999,0,0,0.9500000000000001,0,0,0,0,0,0,0,0,0,0.92,0.5999999999999999,0.92,0.7999999999999998,0,0.9800000000000001,0,1,0.97,0,0.9800000000000001,0,0,0.8000000000001412,0.6000000000002502,0,0.59,0,0)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 11).Value = $kValues[$i]
}
